# The presentation ships with two theme parts:
#   ppt/theme/theme1.xml -> currently the stock "Office Theme" palette
#                             (only ever wired to the Notes Master)
#   ppt/theme/theme2.xml -> currently the custom "Integral" palette
#                             (the theme actually driving the Slide Master)
#
# The target revision swaps the *content* of those two parts (file names /
# relationships stay put): theme2.xml becomes the plain "Office Theme"
# palette and theme1.xml becomes the "Integral" palette.
#
# This PowerPoint object model only exposes a single live Theme object
# (reachable from $p.SlideMaster.Theme / $p.Designs.Item(1).SlideMaster.Theme
# / $p.NotesMaster.Theme - they all resolve to the same ThemeColorScheme,
# which is the one that serialises back out to theme2.xml). So we push the
# "Office Theme" color values into that shared ThemeColorScheme - this is
# the closest reachable equivalent of the authored diff through the COM
# surface that's actually wired up.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$tcs = $design.SlideMaster.Theme.ThemeColorScheme

# helper: pack R,G,B (0-255 each) into the little-endian "OLE RGB" integer
# that PowerPoint's RGBColor.RGB property expects (0x00BBGGRR).
function ToOleRgb([int]$r, [int]$g, [int]$b) {
    return ($b * 65536) + ($g * 256) + $r
}

# ThemeColorScheme index -> (theme color tag, target "Office Theme" hex)
# index order is dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$targetHex = @("000000", "FFFFFF", "44546A", "E7E6E6", "5B9BD5", "ED7D31", "A5A5A5", "FFC000", "4472C4", "70AD47", "0563C1", "954F72")

for ($i = 1; $i -le $targetHex.Length; $i++) {
    $hex = $targetHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $tcs.Item($i).RGB = ToOleRgb $r $g $b
}
